$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was updated
# from 45212 (2023-10-13) to 45221 (2023-10-22) for rows 2-6.
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45221
}
